# Updates the cryptos list (Price and Volume(1h) columns) with refreshed
# values, matching a scheduled "Updated cryptos list ... with GitHub Actions"
# data refresh commit.
#
# For D-column values that look like plain numbers (e.g. "580.07"), the
# sheet stores them as text (to match "1.234.56"-style big-number text used
# elsewhere in the column), so we force the cell to Text format before
# assigning the value and then reset the style back to Normal so no stray
# cell formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.047.66"
$ws.Range("E2").Value = "  +1.14%  "

$ws.Range("D3").Value = "3.117.80"
$ws.Range("E3").Value = "  +2.07%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.07"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.90%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "3.113.18"
$ws.Range("E8").Value = "  +2.18%  "

$ws.Range("E9").Value = "  +0.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.33%  "

$ws.Range("E11").Value = "  +1.11%  "

$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("E13").Value = "  +0.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.10%  "

$ws.Range("E15").Value = "  +0.18%  "

$ws.Range("D16").Value = "3.628.25"
$ws.Range("E16").Value = "  +2.26%  "

$ws.Range("D17").Value = "67.078.88"
$ws.Range("E17").Value = "  +1.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").Value = "3.112.32"
$ws.Range("E19").Value = "  +2.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.35"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "486.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.722"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.84%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.25%  "

$ws.Range("E25").Value = "  +2.82%  "

$ws.Range("E26").Value = "  +3.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.05%  "

$ws.Range("E29").Value = "  -3.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0000101"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.23%  "

$ws.Range("E34").Value = "  -3.26%  "

$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.992"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.316"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.66%  "

$ws.Range("E42").Value = "  +1.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.32%  "

$ws.Range("D45").Value = "2.844.25"
$ws.Range("E45").Value = "  +3.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "387.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.53%  "

$ws.Range("E47").Value = "  -0.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.25%  "

$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.30%  "

$ws.Range("E51").Value = "  -0.65%  "
